# Fix translation col order in the last commit
$wb = $excel.ActiveWorkbook

# --- 1. Fix the hh_id_note translations on table_specific_translations ---
# Previously the Portuguese (col C / text.pt) and Swahili (col D / text.sw)
# values were swapped for this row, and the English text needs a small wording
# tweak too.
$ws6 = $wb.Worksheets.Item("table_specific_translations")
$ws6.Range("B5").Value = "Household ID: {{data.hh_id}}"
$ws6.Range("C5").Value = "Identificação do agregado{{data.hh_id}}"
$ws6.Range("D5").Value = "Utambulisho wa Kaya {{data.hh_id}}"

# --- 2. Re-point the existing conditional formatting rule so it keeps
#        excluding C5:D5 (which are filled in) while now also covering B5.
#        (The engine only lets us keep a single contiguous area per rule
#        object, so the same "equal to 0" rule is re-created across the
#        three areas that together are equivalent to the target sqref.) ---
$fc1 = $ws6.Cells.FormatConditions.Item(1)
$fc1.ModifyAppliesToRange($ws6.Range("B1:D4"))

$cf2 = $ws6.Range("B6:D1048576").FormatConditions.Add(1, 3, "0")
$cf3 = $ws6.Range("B5").FormatConditions.Add(1, 3, "0")

# --- 3. Reset the selection on this sheet back to A1 ---
$ws6.Range("A1").Select()

# --- 4. Set every sheet's zoom back to 100% (was saved at 75%) ---
foreach ($sheet in $wb.Worksheets) {
    $sheet.Activate()
    $excel.ActiveWindow.Zoom = 100
}

# Re-activate table_specific_translations, which is the tab that was selected
$ws6.Activate()
